$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# optimization_parameters sheet: restructure rows (rename "Model" row to
# "production_function", insert a new "L_curve" row after it, drop the
# trailing "Deletion" row, and trim the stray duplicated header cells).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 originally repeats the "value" header across C1:F1 - only A1/B1 stay.
$ws.Range("C1:F1").ClearContents()

# "Model" (row 8) becomes "production_function".
$ws.Range("A8").Value = "production_function"

# Insert the new "L_curve" row right after it.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# The old "Deletion" row (now pushed down to row 17 by the insert above) is
# removed entirely.
$ws.Rows.Item(17).Delete()

# ---------------------------------------------------------------------------
# Sheet-tab / selection bookkeeping: optimization_parameters becomes the
# active tab with C1:F1 selected; production_rates loses the active marker.
# ---------------------------------------------------------------------------
$ws.Select() | Out-Null
$ws.Range("C1:F1").Select() | Out-Null
